$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8810.0625
$ws.Range("J64").Value = 9278
$ws.Range("L64").Value = 9278
$ws.Range("N64").Value = -9774
$ws.Range("H67").Value = 8810.0625
$ws.Range("J67").Value = 9278
$ws.Range("L67").Value = 9278
$ws.Range("N67").Value = -10994
$ws.Range("H74").Value = 6865.276
$ws.Range("I74").Value = 3807.0833
$ws.Range("K74").Value = 3807.0833
$ws.Range("M74").Value = -2871.0833
$ws.Range("H76").Value = 3722.6924
$ws.Range("I76").Value = 3709.6
$ws.Range("K76").Value = 3709.6
$ws.Range("M76").Value = -3394.6
$ws.Range("H77").Value = 6865.276
$ws.Range("I77").Value = 3807.0833
$ws.Range("K77").Value = 19035.4165
$ws.Range("M77").Value = -14355.4165
$ws.Range("H79").Value = 3722.6924
$ws.Range("I79").Value = 3709.6
$ws.Range("K79").Value = 3709.6
$ws.Range("M79").Value = -2617.6
$ws.Range("H99").Value = 1405.2858
$ws.Range("I99").Value = 598
$ws.Range("J99").Value = 2010.75
$ws.Range("K99").Value = 1794
$ws.Range("L99").Value = 6032.25
$ws.Range("M99").Value = -296
$ws.Range("N99").Value = -9028.25
$ws.Range("H112").Value = 10354.407
$ws.Range("I112").Value = 664.3333
$ws.Range("J112").Value = 13123
$ws.Range("K112").Value = 1992.9999
$ws.Range("L112").Value = 39369
$ws.Range("M112").Value = -884.9999
$ws.Range("N112").Value = -41585
$ws.Range("H113").Value = 2796.8572
$ws.Range("J113").Value = 2345
$ws.Range("L113").Value = 2345
$ws.Range("N113").Value = -8853
$ws.Range("H132").Value = 80884.03
$ws.Range("I132").Value = 91697.96000000001
$ws.Range("J132").Value = 10593.5
$ws.Range("K132").Value = 275093.88
$ws.Range("L132").Value = 31780.5
$ws.Range("M132").Value = -272563.88
$ws.Range("N132").Value = -36840.5
$ws.Range("H138").Value = 2171.8809
$ws.Range("I138").Value = 1538.3636
$ws.Range("J138").Value = 2868.75
$ws.Range("K138").Value = 4615.0908
$ws.Range("L138").Value = 8606.25
$ws.Range("M138").Value = 524.9092000000001
$ws.Range("N138").Value = -18886.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 902570.8
$ws.Range("I61").Value = 953889.1
$ws.Range("K61").Value = 953889.1
$ws.Range("M61").Value = -953677.1
$ws.Range("H122").Value = 3142.9092
$ws.Range("I122").Value = 2946.25
$ws.Range("K122").Value = 8838.75
$ws.Range("M122").Value = -6388.75
$ws.Range("H132").Value = 951596.5
$ws.Range("J132").Value = 4881.778
$ws.Range("L132").Value = 14645.334
$ws.Range("N132").Value = -19705.334
$ws.Range("H136").Value = 902570.8
$ws.Range("I136").Value = 953889.1
$ws.Range("K136").Value = 2861667.3
$ws.Range("M136").Value = -2859117.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1024.6
$ws.Range("I20").Value = 849.625
$ws.Range("J20").Value = 1224.5714
$ws.Range("K20").Value = 849.625
$ws.Range("L20").Value = 1224.5714
$ws.Range("M20").Value = -602.625
$ws.Range("N20").Value = -1718.5714
$ws.Range("H44").Value = 39983.5
$ws.Range("J44").Value = 39983.5
$ws.Range("L44").Value = 39983.5
$ws.Range("N44").Value = -40977.5
$ws.Range("H86").Value = 1374.6666
$ws.Range("I86").Value = 1499.5
$ws.Range("J86").Value = 1125
$ws.Range("K86").Value = 1499.5
$ws.Range("L86").Value = 1125
$ws.Range("M86").Value = -376.5
$ws.Range("N86").Value = -3371
$ws.Range("H89").Value = 1374.6666
$ws.Range("I89").Value = 1499.5
$ws.Range("J89").Value = 1125
$ws.Range("K89").Value = 7497.5
$ws.Range("L89").Value = 5625
$ws.Range("M89").Value = -1881.5
$ws.Range("N89").Value = -16857
$ws.Range("H134").Value = 800580.5
$ws.Range("I134").Value = 896368.4399999999
$ws.Range("J134").Value = 513216.78
$ws.Range("K134").Value = 2689105.32
$ws.Range("L134").Value = 1539650.34
$ws.Range("M134").Value = -2686570.32
$ws.Range("N134").Value = -1544720.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 620051.5600000001
$ws.Range("I58").Value = 1030077.75
$ws.Range("K58").Value = 1030077.75
$ws.Range("M58").Value = -1029874.75
$ws.Range("H62").Value = 5082.3335
$ws.Range("J62").Value = 6124.75
$ws.Range("L62").Value = 6124.75
$ws.Range("N62").Value = -7372.75
$ws.Range("H65").Value = 5082.3335
$ws.Range("J65").Value = 6124.75
$ws.Range("L65").Value = 30623.75
$ws.Range("N65").Value = -36863.75
$ws.Range("H98").Value = 90000
$ws.Range("J98").Value = 90000
$ws.Range("L98").Value = 90000
$ws.Range("N98").Value = -94492
$ws.Range("H134").Value = 2349554.5
$ws.Range("I134").Value = 6729.88
$ws.Range("K134").Value = 20189.64
$ws.Range("M134").Value = -17654.64
$ws.Range("H136").Value = 620051.5600000001
$ws.Range("I136").Value = 1030077.75
$ws.Range("K136").Value = 3090233.25
$ws.Range("M136").Value = -3087683.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1023.9
$ws.Range("I34").Value = 177.14285
$ws.Range("J34").Value = 2999.6667
$ws.Range("K34").Value = 531.4285500000001
$ws.Range("L34").Value = 8999.000100000001
$ws.Range("M34").Value = -447.4285500000001
$ws.Range("N34").Value = -9167.000100000001
$ws.Range("H37").Value = 85665.086
$ws.Range("J37").Value = 85665.086
$ws.Range("L37").Value = 256995.258
$ws.Range("N37").Value = -257219.258
$ws.Range("H56").Value = 7125.4
$ws.Range("I56").Value = 7125.4
$ws.Range("K56").Value = 7125.4
$ws.Range("M56").Value = -6595.4
$ws.Range("H121").Value = 16669313
$ws.Range("I121").Value = 50000180
$ws.Range("J121").Value = 3880
$ws.Range("K121").Value = 150000540
$ws.Range("L121").Value = 11640
$ws.Range("M121").Value = -149999230
$ws.Range("N121").Value = -14260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 720.9474
$ws.Range("I2").Value = 800.1429000000001
$ws.Range("K2").Value = 800.1429000000001
$ws.Range("M2").Value = -687.1429000000001
$ws.Range("H70").Value = 7227.857
$ws.Range("I70").Value = 7227.857
$ws.Range("K70").Value = 7227.857
$ws.Range("M70").Value = -6957.857
$ws.Range("H73").Value = 7227.857
$ws.Range("I73").Value = 7227.857
$ws.Range("K73").Value = 7227.857
$ws.Range("M73").Value = -6291.857
$ws.Range("H80").Value = 343994.47
$ws.Range("I80").Value = 429110.66
$ws.Range("K80").Value = 429110.66
$ws.Range("M80").Value = -428112.66
$ws.Range("H83").Value = 343994.47
$ws.Range("I83").Value = 429110.66
$ws.Range("K83").Value = 2145553.3
$ws.Range("M83").Value = -2140561.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 10100.6
$ws.Range("I68").Value = 11125
$ws.Range("K68").Value = 11125
$ws.Range("M68").Value = -10376
$ws.Range("H71").Value = 10100.6
$ws.Range("I71").Value = 11125
$ws.Range("K71").Value = 55625
$ws.Range("M71").Value = -51881
$ws.Range("H122").Value = 5447
$ws.Range("I122").Value = 5185.643
$ws.Range("K122").Value = 15556.929
$ws.Range("M122").Value = -13106.929
